$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 450.16
$ws.Range("I33").Value = 463.21738
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 463.21738
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -234.21738
$ws.Range("N33").Value = -758
$ws.Range("H43").Value = 16925
$ws.Range("I43").Value = 4666
$ws.Range("J43").Value = 26119.25
$ws.Range("K43").Value = 4666
$ws.Range("L43").Value = 26119.25
$ws.Range("M43").Value = -4597
$ws.Range("N43").Value = -26257.25
$ws.Range("H51").Value = 6657.684
$ws.Range("I51").Value = 5316.8335
$ws.Range("J51").Value = 7276.5386
$ws.Range("K51").Value = 5316.8335
$ws.Range("L51").Value = 7276.5386
$ws.Range("M51").Value = -4832.8335
$ws.Range("N51").Value = -8244.5386
$ws.Range("H70").Value = 2685.1843
$ws.Range("I70").Value = 4602.353
$ws.Range("J70").Value = 1133.1904
$ws.Range("K70").Value = 13807.059
$ws.Range("L70").Value = 3399.5712
$ws.Range("M70").Value = -13537.059
$ws.Range("N70").Value = -3939.5712
$ws.Range("H73").Value = 2685.1843
$ws.Range("I73").Value = 4602.353
$ws.Range("J73").Value = 1133.1904
$ws.Range("K73").Value = 13807.059
$ws.Range("L73").Value = 3399.5712
$ws.Range("M73").Value = -12871.059
$ws.Range("N73").Value = -5271.5712
$ws.Range("H76").Value = 16136278
$ws.Range("I76").Value = 38472548
$ws.Range("J76").Value = 4525.6665
$ws.Range("K76").Value = 38472548
$ws.Range("L76").Value = 4525.6665
$ws.Range("M76").Value = -38472233
$ws.Range("N76").Value = -5155.6665
$ws.Range("H79").Value = 16136278
$ws.Range("I79").Value = 38472548
$ws.Range("J79").Value = 4525.6665
$ws.Range("K79").Value = 38472548
$ws.Range("L79").Value = 4525.6665
$ws.Range("M79").Value = -38471456
$ws.Range("N79").Value = -6709.6665
$ws.Range("H129").Value = 1342
$ws.Range("J129").Value = 1702.8572
$ws.Range("L129").Value = 5108.571599999999
$ws.Range("N129").Value = -15108.5716
$ws.Range("H132").Value = 3760.6052
$ws.Range("I132").Value = 921.96875
$ws.Range("K132").Value = 2765.90625
$ws.Range("M132").Value = -235.90625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 360504.12
$ws.Range("I61").Value = 279387.53
$ws.Range("J61").Value = 506514
$ws.Range("K61").Value = 279387.53
$ws.Range("L61").Value = 506514
$ws.Range("M61").Value = -279175.53
$ws.Range("N61").Value = -506938
$ws.Range("H74").Value = 291417.8
$ws.Range("I74").Value = 385708.38
$ws.Range("J74").Value = 87121.586
$ws.Range("K74").Value = 385708.38
$ws.Range("L74").Value = 87121.586
$ws.Range("M74").Value = -384834.38
$ws.Range("N74").Value = -88869.586
$ws.Range("H77").Value = 291417.8
$ws.Range("I77").Value = 385708.38
$ws.Range("J77").Value = 87121.586
$ws.Range("K77").Value = 1928541.9
$ws.Range("L77").Value = 435607.93
$ws.Range("M77").Value = -1924173.9
$ws.Range("N77").Value = -444343.93
$ws.Range("H136").Value = 360504.12
$ws.Range("I136").Value = 279387.53
$ws.Range("J136").Value = 506514
$ws.Range("K136").Value = 838162.5900000001
$ws.Range("L136").Value = 1519542
$ws.Range("M136").Value = -835612.5900000001
$ws.Range("N136").Value = -1524642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 514.86206
$ws.Range("I64").Value = 470
$ws.Range("J64").Value = 614.55554
$ws.Range("K64").Value = 470
$ws.Range("L64").Value = 614.55554
$ws.Range("M64").Value = -245
$ws.Range("N64").Value = -1064.55554
$ws.Range("H67").Value = 514.86206
$ws.Range("I67").Value = 470
$ws.Range("J67").Value = 614.55554
$ws.Range("K67").Value = 470
$ws.Range("L67").Value = 614.55554
$ws.Range("M67").Value = 310
$ws.Range("N67").Value = -2174.55554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 415.75
$ws.Range("I22").Value = 329.85715
$ws.Range("J22").Value = 536
$ws.Range("K22").Value = 329.85715
$ws.Range("L22").Value = 536
$ws.Range("M22").Value = 20.14285000000001
$ws.Range("N22").Value = -1236
$ws.Range("H94").Value = 6834.2666
$ws.Range("I94").Value = 1516.25
$ws.Range("J94").Value = 12912
$ws.Range("K94").Value = 1516.25
$ws.Range("L94").Value = 12912
$ws.Range("M94").Value = -1065.25
$ws.Range("N94").Value = -13814

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1936.8206
$ws.Range("I131").Value = 3654.5
$ws.Range("J131").Value = 1624.5151
$ws.Range("K131").Value = 10963.5
$ws.Range("L131").Value = 4873.5453
$ws.Range("M131").Value = -5923.5
$ws.Range("N131").Value = -14953.5453

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1703.9412
$ws.Range("I16").Value = 1769.0714
$ws.Range("J16").Value = 1400
$ws.Range("K16").Value = 1769.0714
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = -1599.0714
$ws.Range("N16").Value = -1740
$ws.Range("H68").Value = 3689.6
$ws.Range("I68").Value = 3018.3635
$ws.Range("J68").Value = 4510
$ws.Range("K68").Value = 3018.3635
$ws.Range("L68").Value = 4510
$ws.Range("M68").Value = -2269.3635
$ws.Range("N68").Value = -6008
$ws.Range("H71").Value = 3689.6
$ws.Range("I71").Value = 3018.3635
$ws.Range("J71").Value = 4510
$ws.Range("K71").Value = 15091.8175
$ws.Range("L71").Value = 22550
$ws.Range("M71").Value = -11347.8175
$ws.Range("N71").Value = -30038
$ws.Range("H82").Value = 1951.375
$ws.Range("I82").Value = 1575.5
$ws.Range("J82").Value = 2327.25
$ws.Range("K82").Value = 1575.5
$ws.Range("L82").Value = 2327.25
$ws.Range("M82").Value = -1214.5
$ws.Range("N82").Value = -3049.25
$ws.Range("H85").Value = 1951.375
$ws.Range("I85").Value = 1575.5
$ws.Range("J85").Value = 2327.25
$ws.Range("K85").Value = 1575.5
$ws.Range("L85").Value = 2327.25
$ws.Range("M85").Value = -327.5
$ws.Range("N85").Value = -4823.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4435.2354
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 4426.6
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 4426.6
$ws.Range("M62").Value = -3876
$ws.Range("N62").Value = -5674.6
$ws.Range("H65").Value = 4435.2354
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 4426.6
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 22133
$ws.Range("M65").Value = -19380
$ws.Range("N65").Value = -28373
$ws.Range("H100").Value = 6881.909
$ws.Range("I100").Value = 2962.625
$ws.Range("J100").Value = 17333.334
$ws.Range("K100").Value = 5925.25
$ws.Range("L100").Value = 34666.668
$ws.Range("M100").Value = -5384.25
$ws.Range("N100").Value = -35748.668
$ws.Range("H106").Value = 29999
$ws.Range("J106").Value = 29999
$ws.Range("L106").Value = 29999
$ws.Range("N106").Value = -32523
